# Scheduled runner update: refresh market-price-derived columns (H:N)
# on the Leve profit sheets. Values pulled from the latest market data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76: H76,I76,J76,K76,L76,M76,N76
$ws.Range("H76").Value = 3176.923
$ws.Range("I76").Value = 2833.3333
$ws.Range("J76").Value = 3471.4285
$ws.Range("K76").Value = 2833.3333
$ws.Range("L76").Value = 3471.4285
$ws.Range("M76").Value = -2518.3333
$ws.Range("N76").Value = -4101.4285
# Row 79: H79,I79,J79,K79,L79,M79,N79
$ws.Range("H79").Value = 3176.923
$ws.Range("I79").Value = 2833.3333
$ws.Range("J79").Value = 3471.4285
$ws.Range("K79").Value = 2833.3333
$ws.Range("L79").Value = 3471.4285
$ws.Range("M79").Value = -1741.3333
$ws.Range("N79").Value = -5655.4285
# Row 129: H129,J129,L129,N129
$ws.Range("H129").Value = 1462.8276
$ws.Range("J129").Value = 1513.5927
$ws.Range("L129").Value = 4540.7781
$ws.Range("N129").Value = -14540.7781
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 6029.9785
$ws.Range("I132").Value = 4578.028
$ws.Range("K132").Value = 13734.084
$ws.Range("M132").Value = -11204.084
# Row 135: H135,I135,J135,K135,L135,M135,N135
$ws.Range("H135").Value = 5080.154
$ws.Range("I135").Value = 421.7143
$ws.Range("J135").Value = 10515
$ws.Range("K135").Value = 3795.4287
$ws.Range("L135").Value = 94635
$ws.Range("M135").Value = -1260.4287
$ws.Range("N135").Value = -99705

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32,I32,K32,M32
$ws.Range("H32").Value = 3884347.5
$ws.Range("I32").Value = 4007.3472
$ws.Range("K32").Value = 4007.3472
$ws.Range("M32").Value = -3720.3472
# Row 63: H63,I63,J63,K63,L63,M63,N63
$ws.Range("H63").Value = 2328.6365
$ws.Range("I63").Value = 1966.5
$ws.Range("J63").Value = 5950
$ws.Range("K63").Value = 1966.5
$ws.Range("L63").Value = 5950
$ws.Range("M63").Value = -1280.5
$ws.Range("N63").Value = -7322
# Row 66: H66,I66,J66,K66,L66,M66,N66
$ws.Range("H66").Value = 2328.6365
$ws.Range("I66").Value = 1966.5
$ws.Range("J66").Value = 5950
$ws.Range("K66").Value = 9832.5
$ws.Range("L66").Value = 29750
$ws.Range("M66").Value = -6400.5
$ws.Range("N66").Value = -36614
# Row 74: H74,I74,K74,M74
$ws.Range("H74").Value = 14286999
$ws.Range("I74").Value = 20001342
$ws.Range("K74").Value = 20001342
$ws.Range("M74").Value = -20000468
# Row 77: H77,I77,K77,M77
$ws.Range("H77").Value = 14286999
$ws.Range("I77").Value = 20001342
$ws.Range("K77").Value = 100006710
$ws.Range("M77").Value = -100002342
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 4814091.5
$ws.Range("I132").Value = 2622.6667
$ws.Range("J132").Value = 7361339.5
$ws.Range("K132").Value = 7868.000100000001
$ws.Range("L132").Value = 22084018.5
$ws.Range("M132").Value = -5338.000100000001
$ws.Range("N132").Value = -22089078.5

$ws = $wb.Worksheets.Item("BSM")
# Row 105: H105,I105,J105,K105,L105,M105,N105
$ws.Range("H105").Value = 2651.4285
$ws.Range("I105").Value = 1375
$ws.Range("J105").Value = 6736
$ws.Range("K105").Value = 1375
$ws.Range("L105").Value = 6736
$ws.Range("M105").Value = 372
$ws.Range("N105").Value = -10230
# Row 134: H134,I134,J134,K134,L134,M134,N134
$ws.Range("H134").Value = 7292.1724
$ws.Range("I134").Value = 4013
$ws.Range("J134").Value = 8541.380999999999
$ws.Range("K134").Value = 12039
$ws.Range("L134").Value = 25624.143
$ws.Range("M134").Value = -9504
$ws.Range("N134").Value = -30694.143

$ws = $wb.Worksheets.Item("CRP")
# Row 86: H86,I86,J86,K86,L86,M86,N86
$ws.Range("H86").Value = 52636364
$ws.Range("I86").Value = 90912390
$ws.Range("J86").Value = 6826
$ws.Range("K86").Value = 90912390
$ws.Range("L86").Value = 6826
$ws.Range("M86").Value = -90911267
$ws.Range("N86").Value = -9072
# Row 89: H89,I89,J89,K89,L89,M89,N89
$ws.Range("H89").Value = 52636364
$ws.Range("I89").Value = 90912390
$ws.Range("J89").Value = 6826
$ws.Range("K89").Value = 454561950
$ws.Range("L89").Value = 34130
$ws.Range("M89").Value = -454556334
$ws.Range("N89").Value = -45362
# Row 99: H99,I99,J99,K99,L99,M99,N99
$ws.Range("H99").Value = 3331.923
$ws.Range("I99").Value = 1787.7142
$ws.Range("J99").Value = 5133.5
$ws.Range("K99").Value = 1787.7142
$ws.Range("L99").Value = 5133.5
$ws.Range("M99").Value = -289.7141999999999
$ws.Range("N99").Value = -8129.5
# Row 108: H108,I108,J108,K108,L108,M108,N108
$ws.Range("H108").Value = 17482.5
$ws.Range("I108").Value = 4000
$ws.Range("J108").Value = 21976.666
$ws.Range("K108").Value = 4000
$ws.Range("L108").Value = 21976.666
$ws.Range("M108").Value = -160
$ws.Range("N108").Value = -29656.666
# Row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 3331.923
$ws.Range("I126").Value = 1787.7142
$ws.Range("J126").Value = 5133.5
$ws.Range("K126").Value = 5363.142599999999
$ws.Range("L126").Value = 15400.5
$ws.Range("M126").Value = -2893.142599999999
$ws.Range("N126").Value = -20340.5
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 2633565.8
$ws.Range("I132").Value = 3126380.8
$ws.Range("J132").Value = 5218.8335
$ws.Range("K132").Value = 9379142.399999999
$ws.Range("L132").Value = 15656.5005
$ws.Range("M132").Value = -9376612.399999999
$ws.Range("N132").Value = -20716.5005
# Row 134: H134,I134,J134,K134,L134,M134,N134
$ws.Range("H134").Value = 3245.318
$ws.Range("I134").Value = 1612.5
$ws.Range("J134").Value = 6102.75
$ws.Range("K134").Value = 4837.5
$ws.Range("L134").Value = 18308.25
$ws.Range("M134").Value = -2302.5
$ws.Range("N134").Value = -23378.25

$ws = $wb.Worksheets.Item("GSM")
# Row 70: H70,I70,J70,K70,L70,M70,N70
$ws.Range("H70").Value = 4990.5713
$ws.Range("I70").Value = 5236.8
$ws.Range("J70").Value = 4375
$ws.Range("K70").Value = 5236.8
$ws.Range("L70").Value = 4375
$ws.Range("M70").Value = -4966.8
$ws.Range("N70").Value = -4915
# Row 73: H73,I73,J73,K73,L73,M73,N73
$ws.Range("H73").Value = 4990.5713
$ws.Range("I73").Value = 5236.8
$ws.Range("J73").Value = 4375
$ws.Range("K73").Value = 5236.8
$ws.Range("L73").Value = 4375
$ws.Range("M73").Value = -4300.8
$ws.Range("N73").Value = -6247
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 1598.3585
$ws.Range("I132").Value = 1379.4615
$ws.Range("J132").Value = 1809.1482
$ws.Range("K132").Value = 4138.3845
$ws.Range("L132").Value = 5427.444600000001
$ws.Range("M132").Value = -1608.3845
$ws.Range("N132").Value = -10487.4446

$ws = $wb.Worksheets.Item("LTW")
# Row 40: H40,I40,J40,K40,L40,M40,N40
$ws.Range("H40").Value = 111115230
$ws.Range("I40").Value = 500001000
$ws.Range("J40").Value = 5015
$ws.Range("K40").Value = 500001000
$ws.Range("L40").Value = 5015
$ws.Range("M40").Value = -500000864
$ws.Range("N40").Value = -5287
# Row 61: H61,I61,J61,K61,L61,M61,N61
$ws.Range("H61").Value = 3592.2666
$ws.Range("I61").Value = 1526.2727
$ws.Range("J61").Value = 9273.75
$ws.Range("K61").Value = 1526.2727
$ws.Range("L61").Value = 9273.75
$ws.Range("M61").Value = -1324.2727
$ws.Range("N61").Value = -9677.75
# Row 113: H113,I113,J113,K113,L113,M113,N113
$ws.Range("H113").Value = 3592.2666
$ws.Range("I113").Value = 1526.2727
$ws.Range("J113").Value = 9273.75
$ws.Range("K113").Value = 1526.2727
$ws.Range("L113").Value = 9273.75
$ws.Range("M113").Value = 643.7273
$ws.Range("N113").Value = -13613.75
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 24393444
$ws.Range("I132").Value = 33335614
$ws.Range("J132").Value = 5708.909
$ws.Range("K132").Value = 100006842
$ws.Range("L132").Value = 17126.727
$ws.Range("M132").Value = -100004312
$ws.Range("N132").Value = -22186.727

$ws = $wb.Worksheets.Item("WVR")
# Row 39: H39,I39,K39
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
# M39 no longer applicable for this row (NQ/HQ price converged) -> remove cell
$ws.Range("M39").ClearContents()
